$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Part 1: rewrite the "search depth" sentence and relocate the _GoBack
# bookmark into the middle of it.
# ----------------------------------------------------------------------

# 1a. Replace the whole original sentence fragment with the new wording
#     (this inevitably collapses all runs of the paragraph into one, so
#     we restore the run boundaries afterwards).
$d.Content.Find.Execute(
    "The search depth for my algorithm is 15, but can be easily adjusted by changing a ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The search depth for my algorithm is 15 to allow for quick responses by the AI, but this value can be easily adjusted by changing a ",
    2) | Out-Null

# 1b. Helper: force a run split at the point right after $searchText by
#     dropping a temporary bookmark there (bookmarks force a run break)
#     and then removing the bookmark again, leaving the split in place.
function Split-After([string]$searchText) {
    $r = $d.Content
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $false) | Out-Null
    $r.Collapse(0)
    $d.Bookmarks.Add("ZZZTMPSPLITZZZ", $r) | Out-Null
    $d.Bookmarks("ZZZTMPSPLITZZZ").Delete()
}

# 1c. Recreate every run boundary that must exist in the final text:
#     - the four new boundaries introduced by the edit
#     - the boundaries that already existed before "cons"/"t"/" "/
#       "int value in the code." so those runs stay separate
Split-After "The search depth for my algorithm is 15"
Split-After "The search depth for my algorithm is 15 to allow for quick responses by the AI"
Split-After "The search depth for my algorithm is 15 to allow for quick responses by the AI, but "
Split-After "The search depth for my algorithm is 15 to allow for quick responses by the AI, but this value can be easily adjusted by changing a "
Split-After "The search depth for my algorithm is 15 to allow for quick responses by the AI, but this value can be easily adjusted by changing a cons"
Split-After "The search depth for my algorithm is 15 to allow for quick responses by the AI, but this value can be easily adjusted by changing a const"
Split-After "The search depth for my algorithm is 15 to allow for quick responses by the AI, but this value can be easily adjusted by changing a const "
Split-After "The search depth for my algorithm is 15 to allow for quick responses by the AI, but this value can be easily adjusted by changing a const int value in the code."

# 1d. Remove the _GoBack bookmark from its old location (between "...way
#     that it" and " arbitrates...") before re-adding it, so only one
#     instance of the bookmark ever exists at a time.
$d.Bookmarks("_GoBack").Delete()

# 1e. Drop the (permanent) _GoBack bookmark right after "this value ",
#     i.e. just before "can be easily adjusted by changing a ".
$r = $d.Content
$r.Find.Execute(
    "The search depth for my algorithm is 15 to allow for quick responses by the AI, but this value ",
    $true, $false, $false, $false, $false, $true, 1, $false, $false) | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
